# Insert a new weekly Mango price record as row 177, pushing the
# existing rows 177-199 down to 178-200 (dimension grows to A1:T200).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 177; Excel shifts rows 177:199
# down to 178:200 and extends the used range accordingly.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new record.
$ws.Cells.Item(177, 1).Value = 7
$ws.Cells.Item(177, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(177, 3).Value = "Ñuble"
$ws.Cells.Item(177, 4).Value = 45212
$ws.Cells.Item(177, 5).Value = 16
$ws.Cells.Item(177, 6).Value = "Fruta"
$ws.Cells.Item(177, 7).Value = 100108
$ws.Cells.Item(177, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(177, 9).Value = 100108002
$ws.Cells.Item(177, 10).Value = "Mango"
$ws.Cells.Item(177, 11).Value = "Sin especificar"
$ws.Cells.Item(177, 12).Value = "Primera"
$ws.Cells.Item(177, 13).Value = 90
$ws.Cells.Item(177, 14).Value = 10000
$ws.Cells.Item(177, 15).Value = 11000
$ws.Cells.Item(177, 16).Value = 10444
$ws.Cells.Item(177, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(177, 18).Value = "Brasil"
$ws.Cells.Item(177, 19).Value = 2611
$ws.Cells.Item(177, 20).Value = 4
